$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row (2..236).
# All of these values move from 45192 (2023-09-23) to 45202 (2023-10-03).
for ($r = 2; $r -le 236; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}
